# Auto-generated edit script applying the Atomos_Profits diff
# Updates leve profit calculation columns (H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 15625661
$ws.Range("I19").Value = 62500376
$ws.Range("J19").Value = 756.75
$ws.Range("K19").Value = 62500376
$ws.Range("L19").Value = 756.75
$ws.Range("M19").Value = -62500201
$ws.Range("N19").Value = -1106.75

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H87").Value = 29346.4
$ws.Range("J87").Value = 29346.4
$ws.Range("L87").Value = 29346.4
$ws.Range("N87").Value = -31842.4

$ws.Range("H90").Value = 29346.4
$ws.Range("J90").Value = 29346.4
$ws.Range("L90").Value = 88039.20000000001
$ws.Range("N90").Value = -100519.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 50003700
$ws.Range("I2").Value = 83335500
$ws.Range("J2").Value = 6000
$ws.Range("K2").Value = 83335500
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = -83335387
$ws.Range("N2").Value = -6226

$ws.Range("H32").Value = 4620.0244
$ws.Range("I32").Value = 3230.5146
$ws.Range("J32").Value = 11369.071
$ws.Range("K32").Value = 3230.5146
$ws.Range("L32").Value = 11369.071
$ws.Range("M32").Value = -2943.5146
$ws.Range("N32").Value = -11943.071

$ws.Range("H45").Value = 1732.091
$ws.Range("I45").Value = 1088.16
$ws.Range("K45").Value = 1088.16
$ws.Range("M45").Value = -711.1600000000001

$ws.Range("H52").Value = 39780
$ws.Range("J52").Value = 39780
$ws.Range("L52").Value = 39780
$ws.Range("N52").Value = -40416

$ws.Range("H61").Value = 3663.6428
$ws.Range("I61").Value = 1729.1
$ws.Range("J61").Value = 8500
$ws.Range("K61").Value = 1729.1
$ws.Range("L61").Value = 8500
$ws.Range("M61").Value = -1517.1
$ws.Range("N61").Value = -8924

$ws.Range("H116").Value = 50003700
$ws.Range("I116").Value = 83335500
$ws.Range("J116").Value = 6000
$ws.Range("K116").Value = 83335500
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = -83333206
$ws.Range("N116").Value = -10588

$ws.Range("H132").Value = 27030380
$ws.Range("I132").Value = 28574514
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 85723542
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -85721012
$ws.Range("N132").Value = -29060

$ws.Range("H136").Value = 3663.6428
$ws.Range("I136").Value = 1729.1
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 5187.299999999999
$ws.Range("L136").Value = 25500
$ws.Range("M136").Value = -2637.299999999999
$ws.Range("N136").Value = -30600

$ws.Range("H141").Value = 36898.633
$ws.Range("J141").Value = 36898.633
$ws.Range("L141").Value = 36898.633
$ws.Range("N141").Value = -47258.633

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 50003700
$ws.Range("I3").Value = 83335500
$ws.Range("J3").Value = 6000
$ws.Range("K3").Value = 83335500
$ws.Range("L3").Value = 6000
$ws.Range("M3").Value = -83335386
$ws.Range("N3").Value = -6228

$ws.Range("H94").Value = 383.70587
$ws.Range("I94").Value = 374.86667
$ws.Range("J94").Value = 450
$ws.Range("K94").Value = 374.86667
$ws.Range("L94").Value = 450
$ws.Range("M94").Value = 76.13333
$ws.Range("N94").Value = -1352

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4563.1665
$ws.Range("I16").Value = 4126.3335
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 4126.3335
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -3839.3335
$ws.Range("N16").Value = -5574

$ws.Range("H22").Value = 941.4286
$ws.Range("I22").Value = 295
$ws.Range("J22").Value = 3010
$ws.Range("K22").Value = 295
$ws.Range("L22").Value = 3010
$ws.Range("M22").Value = 55
$ws.Range("N22").Value = -3710

$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 10000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -12996

$ws.Range("H113").Value = 4563.1665
$ws.Range("I113").Value = 4126.3335
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4126.3335
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -1956.3335
$ws.Range("N113").Value = -9340

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -34940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1478.1052
$ws.Range("I5").Value = 528.7692
$ws.Range("K5").Value = 1586.3076
$ws.Range("M5").Value = -1474.3076

$ws.Range("H122").Value = 1808.1818
$ws.Range("J122").Value = 1939
$ws.Range("L122").Value = 17451
$ws.Range("N122").Value = -22351

$ws.Range("H131").Value = 976.92725
$ws.Range("J131").Value = 1170.875
$ws.Range("L131").Value = 3512.625
$ws.Range("N131").Value = -13592.625

$ws.Range("H135").Value = 1478.1052
$ws.Range("I135").Value = 528.7692
$ws.Range("K135").Value = 4758.922799999999
$ws.Range("M135").Value = -2223.922799999999

$ws.Range("H141").Value = 2326.8462
$ws.Range("I141").Value = 2326.8462
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6980.5386
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1800.5386
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5018.421
$ws.Range("I122").Value = 9200
$ws.Range("J122").Value = 3903.3333
$ws.Range("K122").Value = 27600
$ws.Range("L122").Value = 11709.9999
$ws.Range("M122").Value = -25150
$ws.Range("N122").Value = -16609.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 66667764
$ws.Range("I22").Value = 166666990
$ws.Range("J22").Value = 1609.3334
$ws.Range("K22").Value = 166666990
$ws.Range("L22").Value = 1609.3334
$ws.Range("M22").Value = -166666695
$ws.Range("N22").Value = -2199.3334

$ws.Range("H27").Value = 66667764
$ws.Range("I27").Value = 166666990
$ws.Range("J27").Value = 1609.3334
$ws.Range("K27").Value = 166666990
$ws.Range("L27").Value = 1609.3334
$ws.Range("M27").Value = -166666883
$ws.Range("N27").Value = -1823.3334

$ws.Range("H61").Value = 200002700
$ws.Range("I61").Value = 250000880
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 250000880
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -250000678
$ws.Range("N61").Value = -10404

$ws.Range("H113").Value = 200002700
$ws.Range("I113").Value = 250000880
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -249998710
$ws.Range("N113").Value = -14340
